$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/16/2025  Through  6/22/2025"

# --- Integer-format cells (#,##0) ---
$intCells = @{
  "C15" = 2
  "F15" = 2
  "G15" = 2
  "I15" = 17
  "J15" = 8
  "D16" = 2
  "F16" = 10
  "G16" = 6
  "I16" = 38
  "J16" = 41
  "C17" = 6
  "D17" = 7
  "F17" = 22
  "G17" = 18
  "I17" = 108
  "J17" = 105
  "C18" = 2
  "D18" = 4
  "F18" = 4
  "G18" = 12
  "I18" = 36
  "J18" = 38
  "C19" = 7
  "D19" = 3
  "F19" = 34
  "G19" = 22
  "I19" = 186
  "J19" = 167
  "C20" = 2
  "D20" = 2
  "F20" = 5
  "G20" = 13
  "I20" = 28
  "J20" = 37
  "C21" = 19
  "D21" = 18
  "F21" = 77
  "G21" = 73
  "I21" = 413
  "J21" = 396
  "C23" = 3
  "D23" = 1
  "F23" = 3
  "G23" = 1
  "I23" = 9
  "J23" = 5
  "C24" = 25
  "D24" = 28
  "F24" = 89
  "G24" = 107
  "I24" = 697
  "J24" = 707
  "C25" = 10
  "D25" = 21
  "F25" = 64
  "G25" = 72
  "I25" = 479
  "J25" = 488
  "C26" = 14
  "D26" = 12
  "F26" = 47
  "G26" = 43
  "I26" = 277
  "J26" = 285
  "C27" = 2
  "F27" = 3
  "G27" = 4
  "I27" = 20
  "J27" = 16
  "C28" = 2
  "D28" = 2
  "F28" = 4
  "G28" = 5
  "I28" = 29
  "J28" = 29
  "I29" = 1
  "J29" = 1
  "I30" = 1
  "J30" = 1
  "G31" = 2
  "J31" = 6
  "D33" = 1
  "G33" = 1
  "I33" = 1
  "J33" = 3
  "J39" = 4
  "J40" = 18
  "J41" = 89
  "J42" = 211
  "J43" = 102
  "J44" = 356
  "J45" = 91
}
foreach ($k in $intCells.Keys) {
  $ws.Range($k).Value = $intCells[$k]
  $ws.Range($k).NumberFormat = "#,##0"
}

# --- Decimal-format cells (1 dp, #,##0.0;-#,##0.0) ---
$decCells = @{
  "L14" = -100
  "H15" = 0
  "K15" = 112.5
  "L15" = 240
  "E16" = -100
  "H16" = 66.666666666666
  "K16" = -7.317073170731
  "L16" = 2.702702702702
  "E17" = -14.285714285714
  "H17" = 22.222222222222
  "K17" = 2.857142857142
  "L17" = -3.571428571428
  "E18" = -50
  "H18" = -66.666666666666
  "K18" = -5.263157894736
  "L18" = 50
  "E19" = 133.333333333333
  "H19" = 54.545454545454
  "K19" = 11.377245508982
  "L19" = -3.125
  "E20" = 0
  "H20" = -61.538461538461
  "K20" = -24.324324324324
  "L20" = -48.148148148148
  "E23" = 200
  "H23" = 200
  "K23" = 80
  "L23" = -10
  "E24" = -10.714285714285
  "H24" = -16.822429906542
  "K24" = -1.414427157001
  "L24" = 8.566978193146
  "E25" = -52.380952380952
  "H25" = -11.111111111111
  "K25" = -1.844262295081
  "L25" = 23.772609819121
  "E26" = 16.666666666666
  "H26" = 9.302325581395
  "K26" = -2.807017543859
  "L26" = 7.364341085271
  "H27" = -25
  "K27" = 25
  "L27" = 122.222222222222
  "E28" = 0
  "H28" = -20
  "K28" = 0
  "L28" = -19.444444444444
  "K29" = 0
  "L29" = -50
  "K30" = 0
  "L30" = -50
  "H31" = -100
  "K31" = -100
  "L31" = -100
  "E33" = -100
  "H33" = -100
  "K33" = -66.666666666666
}
foreach ($k in $decCells.Keys) {
  $ws.Range($k).Value = $decCells[$k]
  $ws.Range($k).NumberFormat = '#,##0.0;"-"#,##0.0'
}

# --- Decimal-format cells (2 dp, row 21 totals, #,##0.00;-#,##0.00) ---
$dec2Cells = @{
  "E21" = 5.555555555555
  "H21" = 5.479452054794
  "K21" = 4.292929292929
  "L21" = -2.823529411764
}
foreach ($k in $dec2Cells.Keys) {
  $ws.Range($k).Value = $dec2Cells[$k]
  $ws.Range($k).NumberFormat = '#,##0.00;"-"#,##0.00'
}

# --- Text placeholder cells ('0' or '***.*') matching existing style (copy format from C14) ---
$donor = $ws.Range("C14")
$ws.Range("C16").Value = "'0"
$donor.Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D31").Value = "'0"
$donor.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "***.*"
$donor.Copy()
$ws.Range("E31").PasteSpecial(-4122)

Write-Host "Edit complete"